$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    [void]$range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

# --- GDP growth heading & paragraph ---
Replace-Text "GDP growth in 2019 slowed down" "GDP growth in 2019 decreased"
Replace-Text "GDP growth slowed down from 2.5% in 2018 to 1.3% in 2019." "GDP growth decreased from 2.5% in 2018 to 1.3% in 2019."
Replace-Text "Gross capital formation added 0.8pp." "Gross capital formation shared 0.8pp."
Replace-Text "Government expenditure added 0.4pp." "Government expenditure gave 0.4pp."
Replace-Text "net exports shaved 1.4pp from growth." "net exports subtracted 1.4pp from growth."
Replace-Text "Industry (including construction) gave 0.3pp." "Industry (including construction) shared 0.3pp."

# --- Gross capital formation (demand side) paragraph ---
Replace-Text "Gross capital formation expanded by the biggest margin at 3.8% annual growth." "Gross capital formation grew by the biggest margin at 3.8% annual growth."
Replace-Text "Private consumption increased by 2.5%." "Private consumption grew by 2.5%."
Replace-Text "Government expenditure jumped by 2.2%." "Government expenditure grew by 2.2%."

# --- Services (supply side) paragraph ---
Replace-Text "Services expanded by the largest edge at 1.6% annual growth." "Services increased by the largest edge at 1.6% annual growth."
Replace-Text "Industry (including construction) expanded by 1.0%." "Industry (including construction) grew by 1.0%."
Replace-Text "Agriculture increased by 0.6%." "Agriculture grew by 0.6%."

# --- Unemployment / inflation heading & paragraph ---
Replace-Text "Unemployment improved; inflation jumped" "Unemployment plunged; inflation increased"
Replace-Text "Unemployment rate improved from 4.8% in 2018 to 4.5% in 2019. Consequently, inflation jumped from 2.9% to 4.5%." "Unemployment rate plunged from 4.8% in 2018 to 4.5% in 2019. Consequently, inflation increased from 2.9% to 4.5%."

# --- Output contracted in Q2 2020 paragraph ---
Replace-Text "Net exports expanded by the biggest margin at 227.5% annual growth." "Net exports picked up by the biggest margin at 227.5% annual growth."
Replace-Text "Government expenditure increased by 1.6%." "Government expenditure expanded by 1.6%."
Replace-Text "private consumption and gross capital formation contracted by 22.2% and 6.2%, respectively." "private consumption and gross capital formation decreased by 22.2% and 6.2%, respectively."

# --- Retail sales paragraph ---
Replace-Text "Growth in the retail sector jumped from a contraction of 3.1% in September" "Growth in the retail sector increased from a contraction of 3.1% in September"

# --- Consumer confidence paragraph ---
Replace-Text "Confidence increased from -30.0 points in the previous quarter." "Confidence improved from -30.0 points in the previous quarter."

# --- Remove "Industrial output shrank" heading and its paragraph ---
$headingIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd("`r","`a") -eq "Industrial output shrank") {
        $headingIdx = $i
        break
    }
}
$headingPara = $d.Paragraphs.Item($headingIdx)
$bodyPara = $d.Paragraphs.Item($headingIdx + 1)
$delRange = $d.Range($headingPara.Range.Start, $bodyPara.Range.End)
$delRange.Delete()

# --- "Inflation worsened" heading -> "Inflation jumped" ---
Replace-Text "Inflation worsened" "Inflation jumped"

# --- Overall inflation paragraph: full rewrite ---
Replace-Text "Overall inflation worsened to 4.0% year-on-year in October from 3.7% in the previous month. Prices for food products worsened to 4.8% from 4.3%, while housing, rent, water, electricity, gas & other fuels jumped to 3.4% from 3.3%. Transportation slowed down to -0.1% from 0.5%, while communication slowed down to 4.3% from 5.2%. Meanwhile, prices for health/medical care worsened to 4.2% from 4.1%, recreation rose to 2.3% from 2.0%, and education slowed down to 1.9% from 2.1%" "Overall inflation jumped to 4.0% year-on-year in October from 3.7% in the previous month. Prices for food products rose to 4.8% from 4.3%, while housing, rent, water, electricity, gas & other fuels worsened to 3.4% from 3.3%. Transportation improved to -0.1% from 0.5%, while communication improved to 4.3% from 5.2%. Meanwhile, prices for health/medical care rose to 4.2% from 4.1%, recreation worsened to 2.3% from 2.0%, and education slowed down to 1.9% from 2.1%"

# --- Outlook paragraph: drop "On the supply side, " and change -3.5% to 3.2% ---
Replace-Text "On the supply side, industrial production is seen to grow by -3.5% and -3.5%." "Industrial production is seen to grow by -3.5% and 3.2%."
